# Updates cryptos list values (Price / Volume(1h) columns, and two
# coin name/link swaps) to match the latest scrape, per commit:
# "Updated cryptos list on Mon May 22 09:22:24 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.033.42'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '1.828.75'
$ws.Range("E3").Value = '  -0.26%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.006'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.52%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '310.86'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.25%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.006'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.34%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4624'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -2.33%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3746'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.63%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07259'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.54%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.8647'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.43%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '19.98'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -2.56%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.07812'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +6.56%  '
$ws.Range("E13").Value = '  -2.50%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.351'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.48%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '6.539'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.38%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '91.98'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.32%  '
$ws.Range("E17").Value = '  -0.20%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000008705'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").Value = '27.150.74'
$ws.Range("E20").Value = '  -1.62%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '14.56'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.52%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.163'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.35%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.57'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.02%  '
$ws.Range("D24").Value = '2.077.47'
$ws.Range("E24").Value = '  -1.41%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '153.39'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.99%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.839'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.88%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.20'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.43%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.097'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.14%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.135'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.88%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '115.43'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -1.57%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.08845'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.63%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.968'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.7291'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.76%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.445'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.24%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.137'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.23%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.492'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.60%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.079'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -1.41%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01945'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.05241'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.85%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.930'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.15%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '7.241'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.06%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.5175'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.23%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.1629'
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.8581'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -15.18%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.212'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.37%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4824'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.01%  '
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("E48").Value = '  -3.46%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '102.82'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.22%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.06260'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("E51").Value = '  -2.74%  '
